$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Range("F44").Value2 = 'Virtus'
$ws.Range("H44").Value2 = 'La Fiorita'
$ws.Range("I44").Value2 = 1
$ws.Range("J44").Value2 = 2.64
$ws.Range("L44").Value2 = 2.82
$ws.Range("N44").Value2 = 2.75
$ws.Range("P44").Value2 = 2.69
$ws.Range("R44").Value2 = 2.47
$ws.Range("T44").Value2 = 2.64
$ws.Range("V44").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/virtus-la-fiorita/Glb5KRXQ/'

# Row 46
$ws.Range("F46").Value2 = 'San Giovanni'
$ws.Range("H46").Value2 = 'Fiorentino'
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 3.06
$ws.Range("L46").Value2 = 3.24
$ws.Range("N46").Value2 = 3.21
$ws.Range("P46").Value2 = 3.24
$ws.Range("R46").Value2 = 1.96
$ws.Range("T46").Value2 = 2.05
$ws.Range("V46").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/san-giovanni-fiorentino/xba1L7IK/'

# Row 48
$ws.Range("F48").Value2 = 'Tre Penne'
$ws.Range("G48").Value2 = 3
$ws.Range("H48").Value2 = 'Murata'
$ws.Range("I48").Value2 = 0
$ws.Range("J48").Value2 = 1.45
$ws.Range("L48").Value2 = 1.85
$ws.Range("M48").Value2 = '29/10/2023 14:24'
$ws.Range("N48").Value2 = 4.07
$ws.Range("P48").Value2 = 3.81
$ws.Range("Q48").Value2 = '29/10/2023 14:32'
$ws.Range("R48").Value2 = 4.66
$ws.Range("T48").Value2 = 3.26
$ws.Range("U48").Value2 = '29/10/2023 14:24'
$ws.Range("V48").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-penne-ss-murata/zZvxF5el/'

# Row 50
$ws.Range("F50").Value2 = 'Cailungo'
$ws.Range("G50").Value2 = 4
$ws.Range("H50").Value2 = 'Faetano'
$ws.Range("I50").Value2 = 1
$ws.Range("J50").Value2 = 2.69
$ws.Range("L50").Value2 = 3.01
$ws.Range("M50").Value2 = '29/10/2023 14:54'
$ws.Range("N50").Value2 = 3.3
$ws.Range("P50").Value2 = 3.7
$ws.Range("Q50").Value2 = '29/10/2023 14:54'
$ws.Range("R50").Value2 = 2.12
$ws.Range("T50").Value2 = 1.99
$ws.Range("U50").Value2 = '29/10/2023 14:54'
$ws.Range("V50").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/cailungo-sc-faetano/MVzYFotr/'

# Row 56
$ws.Range("F56").Value2 = 'Tre Fiori'
$ws.Range("H56").Value2 = 'Juvenes/Dogana'
$ws.Range("I56").Value2 = 1
$ws.Range("J56").Value2 = 1.44
$ws.Range("L56").Value2 = 1.56
$ws.Range("M56").Value2 = '05/11/2023 14:08'
$ws.Range("N56").Value2 = 3.9
$ws.Range("P56").Value2 = 4.05
$ws.Range("Q56").Value2 = '05/11/2023 14:08'
$ws.Range("R56").Value2 = 5.01
$ws.Range("T56").Value2 = 4.62
$ws.Range("U56").Value2 = '05/11/2023 14:08'
$ws.Range("V56").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-juvenes-dogana/M91DEd15/'

# Row 57
$ws.Range("F57").Value2 = 'Murata'
$ws.Range("H57").Value2 = 'Pennarossa'
$ws.Range("I57").Value2 = 0
$ws.Range("J57").Value2 = 1.37
$ws.Range("L57").Value2 = 1.39
$ws.Range("M57").Value2 = '05/11/2023 14:06'
$ws.Range("N57").Value2 = 4.38
$ws.Range("P57").Value2 = 4.73
$ws.Range("Q57").Value2 = '05/11/2023 14:06'
$ws.Range("R57").Value2 = 5.5
$ws.Range("T57").Value2 = 5.7
$ws.Range("U57").Value2 = '05/11/2023 14:06'
$ws.Range("V57").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ss-murata-ss-pennarossa/4YjIDGGB/'

# Row 63
$ws.Range("F63").Value2 = 'Libertas'
$ws.Range("G63").Value2 = 1
$ws.Range("H63").Value2 = 'San Giovanni'
$ws.Range("J63").Value2 = 1.5
$ws.Range("L63").Value2 = 1.47
$ws.Range("M63").Value2 = '12/11/2023 14:06'
$ws.Range("N63").Value2 = 3.96
$ws.Range("P63").Value2 = 4.35
$ws.Range("Q63").Value2 = '12/11/2023 14:48'
$ws.Range("R63").Value2 = 4.89
$ws.Range("T63").Value2 = 5.13
$ws.Range("U63").Value2 = '12/11/2023 14:06'
$ws.Range("V63").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ac-libertas-san-giovanni/IyHiSFV4/'

# Row 65
$ws.Range("F65").Value2 = 'Virtus'
$ws.Range("G65").Value2 = 6
$ws.Range("H65").Value2 = 'Faetano'
$ws.Range("J65").Value2 = 1.21
$ws.Range("L65").Value2 = 1.19
$ws.Range("M65").Value2 = '12/11/2023 14:53'
$ws.Range("N65").Value2 = 5.7
$ws.Range("P65").Value2 = 6.22
$ws.Range("Q65").Value2 = '12/11/2023 14:53'
$ws.Range("R65").Value2 = 8.26
$ws.Range("T65").Value2 = 9.65
$ws.Range("U65").Value2 = '12/11/2023 14:53'
$ws.Range("V65").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/virtus-sc-faetano/lMhUAE0U/'

# New row 73: copy formatting from row 72, then set values
$ws.Range("A72:V72").Copy()
$ws.Range("A73:V73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A73").Value2 = 72
$ws.Range("B73").Value2 = 'san-marino'
$ws.Range("C73").Value2 = 'campionato-sammarinese'
$ws.Range("D73").Value2 = '2023-2024'
$ws.Range("E73").Value2 = 45262.76041666666
$ws.Range("F73").Value2 = 'Tre Fiori'
$ws.Range("G73").Value2 = 2
$ws.Range("H73").Value2 = 'Cailungo'
$ws.Range("I73").Value2 = 1
$ws.Range("J73").Value2 = 1.2
$ws.Range("K73").Value2 = '02/12/2023 07:42'
$ws.Range("L73").Value2 = 1.16
$ws.Range("M73").Value2 = '02/12/2023 17:39'
$ws.Range("N73").Value2 = 5.91
$ws.Range("O73").Value2 = '02/12/2023 07:42'
$ws.Range("P73").Value2 = 6.32
$ws.Range("Q73").Value2 = '02/12/2023 17:40'
$ws.Range("R73").Value2 = 9.09
$ws.Range("S73").Value2 = '02/12/2023 07:42'
$ws.Range("T73").Value2 = 11.92
$ws.Range("U73").Value2 = '02/12/2023 17:40'
$ws.Range("V73").Value2 = 'https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-cailungo/WMkYvgFA/'

Write-Output "done"